$d = $word.ActiveDocument

# Locate the paragraph that holds the existing "Docente(s)" list entry
# (search on a plain-ASCII substring to avoid the PS string round-trip
# mangling accented characters when read back).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*3444370*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph (3444370 - Rita de Cassia ...)"
}

# New names to add, in the order they should appear before the
# existing "3444370 - ..." entry.
$newNames = @(
    "427823 - Adriane Maria Ferreira Milagres",
    "4873328 - Fernando Segato",
    "6007846 - Júlio César dos Santos"
)

# Insert one blank paragraph (inherits the ListBullet style automatically)
# per new name, directly above the anchor paragraph.
$r = $d.Paragraphs.Item($anchorIndex).Range
$r.Collapse(1)
for ($k = 0; $k -lt $newNames.Count; $k++) {
    $r.InsertParagraphBefore()
}

# Fill every newly created (still separate) paragraph with
# "Name" + line-break char.
for ($k = 0; $k -lt $newNames.Count; $k++) {
    $p = $d.Paragraphs.Item($anchorIndex + $k)
    $p.Range.Text = $newNames[$k] + [char]11
}

# Merge the new paragraphs into the anchor paragraph by repeatedly
# deleting the paragraph mark right after the (growing) anchor
# paragraph. This keeps each line as its own <w:r> run (matching the
# existing "LOT2013.../LOT2017..." run-per-line pattern) instead of
# collapsing everything into one run with multiple <w:t>/<w:br/>
# children.
for ($k = 0; $k -lt $newNames.Count; $k++) {
    $p = $d.Paragraphs.Item($anchorIndex)
    $m = $p.Range
    $m.Collapse(0)
    $m.MoveEnd(1, 1)
    $m.Delete()
}
